$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E (Password) currently holds numeric 0 for every data row (2-22).
# Replace with the text "empty" so the cell becomes a shared-string entry
# (matches the existing "empty" string already used in columns F/G).
for ($row = 2; $row -le 22; $row++) {
    $ws.Cells.Item($row, 5).Value = "empty"
}

# Column H (MailCapacity) is new for rows 2-22: every data row gets the
# value 5.
for ($row = 2; $row -le 22; $row++) {
    $ws.Cells.Item($row, 8).Value = 5
}

# Update the active selection on the sheet from J20 to N9.
$ws.Range("N9").Select() | Out-Null
